{"js": "// Replace each three-digit x one-digit multiplication equation with its\n// updated version. The worksheet table holds 25 equation cells; every one\n// of them is swapped for a new equation (same \"N\u00d7N=N\" shape), so we do a\n// literal, case-sensitive, whole-text search & replace per pair. Each\n// \"before\" string is unique in the document, so there is no ambiguity.\nconst replacements = [\n  [\"956\u00d75=4780\", \"758\u00d72=1516\"],\n  [\"216\u00d73=648\", \"256\u00d73=768\"],\n  [\"615\u00d78=4920\", \"679\u00d73=2037\"],\n  [\"356\u00d74=1424\", \"298\u00d75=1490\"],\n  [\"258\u00d75=1290\", \"494\u00d75=2470\"],\n  [\"171\u00d73=513\", \"369\u00d78=2952\"],\n  [\"340\u00d75=1700\", \"247\u00d75=1235\"],\n  [\"326\u00d74=1304\", \"836\u00d74=3344\"],\n  [\"129\u00d73=387\", \"792\u00d78=6336\"],\n  [\"658\u00d75=3290\", \"467\u00d74=1868\"],\n  [\"593\u00d73=1779\", \"283\u00d75=1415\"],\n  [\"244\u00d79=2196\", \"842\u00d79=7578\"],\n  [\"520\u00d73=1560\", \"710\u00d75=3550\"],\n  [\"952\u00d73=2856\", \"803\u00d73=2409\"],\n  [\"549\u00d79=4941\", \"678\u00d77=4746\"],\n  [\"880\u00d79=7920\", \"390\u00d77=2730\"],\n  [\"168\u00d78=1344\", \"360\u00d76=2160\"],\n  [\"155\u00d73=465\", \"278\u00d74=1112\"],\n  [\"797\u00d75=3985\", \"386\u00d76=2316\"],\n  [\"791\u00d79=7119\", \"502\u00d78=4016\"],\n  [\"202\u00d78=1616\", \"847\u00d75=4235\"],\n  [\"464\u00d74=1856\", \"507\u00d74=2028\"],\n  [\"255\u00d75=1275\", \"455\u00d78=3640\"],\n  [\"919\u00d78=7352\", \"326\u00d76=1956\"],\n  [\"343\u00d74=1372\", \"237\u00d73=711\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const found = body.search(before, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit x one-digit multiplication equation with its\n# updated version. The worksheet table holds 25 equation cells; every one\n# of them is swapped for a new equation (same \"N\u00d7N=N\" shape). Each \"before\"\n# string is unique in the document, so a simple Find/Replace per pair,\n# scoped to a single hit, is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"956\u00d75=4780\", \"758\u00d72=1516\"),\n    @(\"216\u00d73=648\", \"256\u00d73=768\"),\n    @(\"615\u00d78=4920\", \"679\u00d73=2037\"),\n    @(\"356\u00d74=1424\", \"298\u00d75=1490\"),\n    @(\"258\u00d75=1290\", \"494\u00d75=2470\"),\n    @(\"171\u00d73=513\", \"369\u00d78=2952\"),\n    @(\"340\u00d75=1700\", \"247\u00d75=1235\"),\n    @(\"326\u00d74=1304\", \"836\u00d74=3344\"),\n    @(\"129\u00d73=387\", \"792\u00d78=6336\"),\n    @(\"658\u00d75=3290\", \"467\u00d74=1868\"),\n    @(\"593\u00d73=1779\", \"283\u00d75=1415\"),\n    @(\"244\u00d79=2196\", \"842\u00d79=7578\"),\n    @(\"520\u00d73=1560\", \"710\u00d75=3550\"),\n    @(\"952\u00d73=2856\", \"803\u00d73=2409\"),\n    @(\"549\u00d79=4941\", \"678\u00d77=4746\"),\n    @(\"880\u00d79=7920\", \"390\u00d77=2730\"),\n    @(\"168\u00d78=1344\", \"360\u00d76=2160\"),\n    @(\"155\u00d73=465\", \"278\u00d74=1112\"),\n    @(\"797\u00d75=3985\", \"386\u00d76=2316\"),\n    @(\"791\u00d79=7119\", \"502\u00d78=4016\"),\n    @(\"202\u00d78=1616\", \"847\u00d75=4235\"),\n    @(\"464\u00d74=1856\", \"507\u00d74=2028\"),\n    @(\"255\u00d75=1275\", \"455\u00d78=3640\"),\n    @(\"919\u00d78=7352\", \"326\u00d76=1956\"),\n    @(\"343\u00d74=1372\", \"237\u00d73=711\")\n)\n\nforeach ($pair in $pairs) {\n    $before = $pair[0]\n    $after = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $before\n    $find.Replacement.Text = $after\n    $find.Forward = $true\n    $find.Wrap = 0\n\n    # wdFindContinue(1) restricted to this exact text; wdReplaceOne(1) so only\n    # the single known occurrence of each \"before\" string is touched.\n    $find.Execute($before, $false, $false, $false, $false, $false, $true, 1, $false, $after, 1)\n}\n"}
